$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) After "Facilité de transmission des informations entre utilisateurs",
#    insert a new list item "Maintenabilité améliorée", then turn the
#    (pre-existing) following empty paragraph into the new home of the
#    "_GoBack" bookmark (which Word will automatically relocate away from
#    the end of the document, since a document can only have one).
# ---------------------------------------------------------------------------
$rngFind = $d.Content
$rngFind.Find.Execute("Facilité de transmission des informations entre utilisateurs",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraFacilite = $rngFind.Paragraphs(1)

$paraFacilite.Range.InsertParagraphAfter()
$paraMaintenabilite = $paraFacilite.Next()
$paraMaintenabilite.Range.Text = "Maintenabilité améliorée"

$paraBlank = $paraMaintenabilite.Next()
$blankStart = $paraBlank.Range.Start

# Insert a temporary character so the bookmark range is non-collapsed
# (collapsed ranges don't anchor new bookmarks reliably), then strip the
# temporary character back out once the bookmark has been anchored.
$tmpRange = $d.Range($blankStart, $blankStart)
$tmpRange.InsertAfter("X")
$bmRange = $d.Range($blankStart, $blankStart + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$clearRange = $d.Range($blankStart, $blankStart + 1)
$clearRange.Text = ""

# ---------------------------------------------------------------------------
# 2) Drop every stray <w:lastRenderedPageBreak/> left over from the last time
#    the document was paginated in Word. Each one is removed by rebuilding
#    its containing paragraph (same paragraph attributes + run content,
#    minus the page-break marker) through InsertXML, which is immune to the
#    run/run-boundary reshuffling that a narrower, mid-paragraph replace
#    would trigger.
# ---------------------------------------------------------------------------
function Remove-LastRenderedPageBreakParagraph {
    param(
        [string]$anchorText,
        [string]$paragraphOpenTag,
        [string]$innerXml
    )

    $rng = $d.Content
    $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $para = $rng.Paragraphs(1)
    $target = $para.Range

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $paragraphOpenTag + $innerXml + '</w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
}

Remove-LastRenderedPageBreakParagraph `
    -anchorText "serveur web. Le serveur web devra donc avoir" `
    -paragraphOpenTag '<w:p w14:paraId="2EDDF220" w14:textId="77777777" w:rsidR="00371379" w:rsidRDefault="00275712" w:rsidP="00860517">' `
    -innerXml ('<w:r><w:t xml:space="preserve">L&#8217;outil sera h&#233;berg&#233; sur un serveur web. Le serveur web devra donc avoir d&#8217;install&#233; et de correctement configur&#233; : </w:t></w:r>' +
               '<w:r><w:br/></w:r>' +
               '<w:r><w:br/><w:t xml:space="preserve">Apache 2.4, MySQL 5.6, </w:t></w:r>' +
               '<w:proofErr w:type="spellStart"/><w:r><w:t>phpMyAdmin</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
               '<w:r><w:t xml:space="preserve"> 4.2.7, serveur mail fonctionnel.</w:t></w:r>')

Remove-LastRenderedPageBreakParagraph `
    -anchorText "Administrateur pourra créer un compte STF" `
    -paragraphOpenTag '<w:p w14:paraId="2E8ACDFD" w14:textId="77777777" w:rsidR="00FE17E3" w:rsidRDefault="00FE17E3" w:rsidP="00FE17E3">' `
    -innerXml '<w:r><w:t>Administrateur pourra créer un compte STF</w:t></w:r>'

Remove-LastRenderedPageBreakParagraph `
    -anchorText "Automatisée" `
    -paragraphOpenTag '<w:p w14:paraId="42926642" w14:textId="38E84C6F" w:rsidR="005C3BAA" w:rsidRDefault="005C3BAA" w:rsidP="005C3BAA"><w:pPr><w:pStyle w:val="Heading3"/></w:pPr>' `
    -innerXml '<w:r><w:tab/><w:t>Automatisée</w:t></w:r>'

Remove-LastRenderedPageBreakParagraph `
    -anchorText "Le nom du matériel serait cliquable" `
    -paragraphOpenTag '<w:p w14:paraId="5A16C769" w14:textId="77777777" w:rsidR="00640502" w:rsidRDefault="00142CB9" w:rsidP="008178ED">' `
    -innerXml ('<w:r><w:t>Le nom du matériel serait cliquable, et renverrai vers la page « matériel »  détaillant les caractéristiques de ce matériel.</w:t></w:r>' +
               '<w:r><w:br/></w:r>')

# ---------------------------------------------------------------------------
# 3) Merge the two runs around the removed page-break in the "Administrateur"
#    paragraph into a single run.
# ---------------------------------------------------------------------------
Remove-LastRenderedPageBreakParagraph `
    -anchorText "Il sera celui en charge de la maintenance" `
    -paragraphOpenTag '<w:p w14:paraId="2B19A2A9" w14:textId="77777777" w:rsidR="006914CE" w:rsidRDefault="006914CE" w:rsidP="006914CE">' `
    -innerXml ('<w:r><w:t>Il sera celui en charge de la maintenance et de la veille sur l&#8217;outil. Il aura accès au mode développeur et pourra se charger de la gestion des utilisateurs, de la gestion du mode développeur.</w:t></w:r>' +
               '<w:r w:rsidR="00053F57"><w:t xml:space="preserve"> Il fera r</w:t></w:r>' +
               '<w:r w:rsidR="00DE1394"><w:t>emonter toutes les anomalies/bogues</w:t></w:r>' +
               '<w:r w:rsidR="00053F57"><w:t xml:space="preserve"> à l&#8217;équipe de développement et sera leur interlocuteur privilégié.</w:t></w:r>' +
               '<w:r><w:br/></w:r>')
